# BB_Build.xlsx edit: insert "Icon" and "ShowName" columns (G, H) before the
# existing "Desc" column (which shifts from G to I).
#
# Icon (G)       = last path segment of the Prefab value (column D)
# ShowName (H)   = a copy of the existing Desc value
# Desc (I)       = the existing Desc value, unchanged (just relocated)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 10

# Move the existing "Desc" column (G) out to I first, then populate G/H.
for ($r = 1; $r -le $lastRow; $r++) {
    $descCell = $ws.Cells.Item($r, 7)
    $descValue = $descCell.Value()

    if ($r -eq 1) {
        $ws.Cells.Item($r, 9).Value = $descValue
        $ws.Cells.Item($r, 7).Value = "Icon"
        $ws.Cells.Item($r, 8).Value = "ShowName"
    } else {
        $prefab = $ws.Cells.Item($r, 4).Value()
        $parts = $prefab.Split("/")
        $icon = $parts[$parts.Length - 1]

        $ws.Cells.Item($r, 9).Value = $descValue
        $ws.Cells.Item($r, 7).Value = $icon
        $ws.Cells.Item($r, 8).Value = $descValue

        $ws.Cells.Item($r, 8).NumberFormat = "@"
        $ws.Cells.Item($r, 9).NumberFormat = "@"
    }
}

# Columns G:I all end up width "11" (character units) in the committed
# file; 10.29 is the COM ColumnWidth value that round-trips to that.
$ws.Range("G1:I1").EntireColumn.ColumnWidth = 10.29

$ws.Range("H10").Select()
